$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.1800000000005
$ws.Range("H2").Value = [double]"3.573784309984784e-05"
$ws.Range("I2").Value = [double]"3.573784309984784e-05"
$ws.Range("L2").Value = 48.47784524283492
$ws.Range("M2").Value = "[24.775595532183303, 72.18009495348655]"
$ws.Range("N2").Value = 0.0001605449539556769
$ws.Range("O2").Value = 0.0001605449539556769
$ws.Range("P2").Value = 1.415131825941348
$ws.Range("Q2").Value = "[0.8365001460008861, 1.9937635058818097]"
$ws.Range("R2").Value = [double]"1.177259725659852e-05"
$ws.Range("S2").Value = [double]"1.177259725659852e-05"
$ws.Range("T2").Value = 60.10512053761784
$ws.Range("U2").Value = "[46.27852204562325, 73.93171902961242]"
$ws.Range("V2").Value = [double]"2.815681021672845e-11"
$ws.Range("W2").Value = [double]"2.815681021672845e-11"
$ws.Range("X2").Value = 19.50882882882922
$ws.Range("Y2").Value = 17.18994994995029
$ws.Range("Z2").Value = 21.82770770770814
$ws.Range("F3").Value = 25.1800000000005
$ws.Range("H3").Value = 0.002070520235338602
$ws.Range("I3").Value = 0.002070520235338602
$ws.Range("L3").Value = 39.69344789302975
$ws.Range("M3").Value = "[13.377361102640506, 66.00953468341899]"
$ws.Range("N3").Value = 0.003956517256510139
$ws.Range("O3").Value = 0.003956517256510139
$ws.Range("P3").Value = 1.289342330302118
$ws.Range("Q3").Value = "[0.5094474573388847, 2.0692372032653505]"
$ws.Range("R3").Value = 0.001741641781479508
$ws.Range("S3").Value = 0.001741641781479508
$ws.Range("T3").Value = 52.5890995448656
$ws.Range("U3").Value = "[37.37180967567073, 67.80638941406046]"
$ws.Range("V3").Value = [double]"1.166163943189247e-08"
$ws.Range("W3").Value = [double]"1.166163943189247e-08"
$ws.Range("X3").Value = 20.01293293293332
$ws.Range("Y3").Value = 16.88748748748781
$ws.Range("Z3").Value = 23.13837837837884
$ws.Range("F4").Value = 25.1800000000005
$ws.Range("H4").Value = [double]"6.467752411420946e-08"
$ws.Range("I4").Value = [double]"6.467752411420946e-08"
$ws.Range("L4").Value = 70.60621771017875
$ws.Range("M4").Value = "[42.83598287876643, 98.37645254159108]"
$ws.Range("N4").Value = [double]"6.144433793098258e-06"
$ws.Range("O4").Value = [double]"6.144433793098258e-06"
$ws.Range("P4").Value = 0.748447499053424
$ws.Range("Q4").Value = "[0.37107901213573147, 1.1258159859711165]"
$ws.Range("R4").Value = 0.0002371930608366934
$ws.Range("S4").Value = 0.0002371930608366934
$ws.Range("T4").Value = 75.56160134981337
$ws.Range("U4").Value = "[61.36627328145536, 89.75692941817138]"
$ws.Range("V4").Value = [double]"5.617728504603292e-14"
$ws.Range("W4").Value = [double]"5.617728504603292e-14"
$ws.Range("X4").Value = 22.18058058058102
$ws.Range("Y4").Value = 20.66826826826868
$ws.Range("Z4").Value = 23.69289289289336
$ws.Range("F5").Value = 25.1800000000005
$ws.Range("H5").Value = [double]"5.389167910818671e-07"
$ws.Range("I5").Value = [double]"5.389167910818671e-07"
$ws.Range("L5").Value = 72.05076866433055
$ws.Range("M5").Value = "[41.72120798443669, 102.3803293442244]"
$ws.Range("N5").Value = [double]"1.877267778205294e-05"
$ws.Range("O5").Value = [double]"1.877267778205294e-05"
$ws.Range("P5").Value = 0.5471843060306538
$ws.Range("Q5").Value = "[0.10692107129334438, 0.9874475407679633]"
$ws.Range("R5").Value = 0.01600319140295392
$ws.Range("S5").Value = 0.01600319140295392
$ws.Range("T5").Value = 68.12551167362334
$ws.Range("U5").Value = "[52.30048431035406, 83.95053903689262]"
$ws.Range("V5").Value = [double]"3.720379559979392e-11"
$ws.Range("W5").Value = [double]"3.720379559979392e-11"
$ws.Range("X5").Value = 22.9871471471476
$ws.Range("Y5").Value = 21.2227827827832
$ws.Range("Z5").Value = 24.75151151151201
$ws.Range("B6").Value = 0
$ws.Range("F6").Value = 23.08000000000017
$ws.Range("H6").Value = 0.0001206579660029128
$ws.Range("I6").Value = 0.0001206579660029128
$ws.Range("L6").Value = 42.48530303417934
$ws.Range("M6").Value = "[20.611683167852476, 64.3589229005062]"
$ws.Range("N6").Value = 0.0003063770390481224
$ws.Range("O6").Value = 0.0003063770390481224
$ws.Range("P6").Value = 0.1320789704211922
$ws.Range("Q6").Value = "[-0.4842895582110396, 0.748447499053424]"
$ws.Range("R6").Value = 0.6680976392579638
$ws.Range("S6").Value = 0.6680976392579638
$ws.Range("T6").Value = 57.54784241236789
$ws.Range("U6").Value = "[44.493689681398635, 70.60199514333715]"
$ws.Range("V6").Value = [double]"1.879674194071868e-11"
$ws.Range("W6").Value = [double]"1.879674194071868e-11"
$ws.Range("X6").Value = 22.594834834835
$ws.Range("Y6").Value = 20.33073073073088
$ws.Range("Z6").Value = 24.85893893893912
$ws.Range("F7").Value = 23.08000000000017
$ws.Range("H7").Value = [double]"1.54127873885912e-06"
$ws.Range("I7").Value = [double]"1.54127873885912e-06"
$ws.Range("L7").Value = 65.81491738830421
$ws.Range("M7").Value = "[35.988134434435125, 95.6417003421733]"
$ws.Range("N7").Value = [double]"5.699992462471215e-05"
$ws.Range("O7").Value = [double]"5.699992462471215e-05"
$ws.Range("P7").Value = 0.4465527095192696
$ws.Range("Q7").Value = "[-0.04402632347373192, 0.9371317425122712]"
$ws.Range("R7").Value = 0.07337009535387451
$ws.Range("S7").Value = 0.07337009535387451
$ws.Range("T7").Value = 72.19709739601922
$ws.Range("U7").Value = "[56.539032240258855, 87.85516255177959]"
$ws.Range("V7").Value = [double]"5.011102643948107e-12"
$ws.Range("W7").Value = [double]"5.011102643948107e-12"
$ws.Range("X7").Value = 21.43967967967984
$ws.Range("Y7").Value = 19.63763763763778
$ws.Range("Z7").Value = 23.24172172172189
$ws.Range("F8").Value = 23.08000000000017
$ws.Range("H8").Value = [double]"1.220658593981128e-05"
$ws.Range("I8").Value = [double]"1.220658593981128e-05"
$ws.Range("L8").Value = 52.04016951974446
$ws.Range("M8").Value = "[26.066566633515023, 78.0137724059739]"
$ws.Range("N8").Value = 0.0002088916856473055
$ws.Range("O8").Value = 0.0002088916856473055
$ws.Range("P8").Value = 0.6226580034141929
$ws.Range("Q8").Value = "[0.1320789704211922, 1.1132370364071935]"
$ws.Range("R8").Value = 0.0140227501178205
$ws.Range("S8").Value = 0.0140227501178205
$ws.Range("T8").Value = 52.7177873256778
$ws.Range("U8").Value = "[38.992733131901645, 66.44284151945396]"
$ws.Range("V8").Value = [double]"8.368596926544569e-10"
$ws.Range("W8").Value = [double]"8.368596926544569e-10"
$ws.Range("X8").Value = 20.79279279279294
$ws.Range("Y8").Value = 18.99075075075089
$ws.Range("Z8").Value = 22.594834834835

Write-Host "done"